# Update the crypto price table cells to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds decimal-looking text (e.g. "29.095.12", "0.9989") that Excel
# would otherwise auto-convert to a number. Force the whole column to Text
# first, write the values, then restore the default "Normal" style so the
# saved file keeps the original (unstyled) cell formatting.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# --- Rows 2-45: price/volume refresh (no coin re-ordering in this block) ---
$ws.Range("D2").Value = '29.095.12'
$ws.Range("D3").Value = '1.831.94'
$ws.Range("D4").Value = '0.9989'
$ws.Range("D5").Value = '239.04'
$ws.Range("D6").Value = '0.6643'
$ws.Range("D7").Value = '0.9999'
$ws.Range("D9").Value = '0.07323'
$ws.Range("D10").Value = '22.68'
$ws.Range("D12").Value = '1.835.32'
$ws.Range("D13").Value = '5.017'
$ws.Range("D15").Value = '85.95'
$ws.Range("D16").Value = '6.128'
$ws.Range("D17").Value = '29.089.73'
$ws.Range("D18").Value = '0.000008193'
$ws.Range("D19").Value = '227.19'
$ws.Range("D21").Value = '0.9995'
$ws.Range("D22").Value = '7.252'
$ws.Range("D23").Value = '0.9997'
$ws.Range("D24").Value = '160.67'
$ws.Range("D25").Value = '0.1421'
$ws.Range("D26").Value = '8.646'
$ws.Range("D27").Value = '17.94'
$ws.Range("D28").Value = '1.494'
$ws.Range("D29").Value = '4.221'
$ws.Range("D30").Value = '4.098'
$ws.Range("D32").Value = '0.05325'
$ws.Range("D33").Value = '1.849'
$ws.Range("D34").Value = '0.7447'
$ws.Range("D37").Value = '1.298.18'
$ws.Range("D39").Value = '2.706'
$ws.Range("D41").Value = '6.031'
$ws.Range("D43").Value = '103.55'
$ws.Range("D44").Value = '1.984.13'

$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  -2.40%  '
$ws.Range("E6").Value = '  -4.38%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -4.70%  '
$ws.Range("E10").Value = '  -3.90%  '
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("E13").Value = '  -2.63%  '
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("E15").Value = '  -5.75%  '
$ws.Range("E16").Value = '  -3.22%  '
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("E19").Value = '  -4.64%  '
$ws.Range("E20").Value = '  -2.02%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  -4.96%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("E25").Value = '  -4.85%  '
$ws.Range("E26").Value = '  -2.86%  '
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("E28").Value = '  -2.51%  '
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("E31").Value = '  -1.17%  '
$ws.Range("E32").Value = '  +4.12%  '
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("E34").Value = '  -3.54%  '
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("E37").Value = '  -2.37%  '
$ws.Range("E38").Value = '  -3.32%  '
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("E41").Value = '  +2.90%  '
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("E43").Value = '  -2.08%  '
$ws.Range("E45").Value = '  -0.88%  '

# --- Rows 46-51: a new coin (BabyDogeCoin) was inserted at rank 44, pushing
# the remaining coins down by one row; Aptos drops off the bottom of the list ---
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000121'
$ws.Range("E46").Value = '  -3.57%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '63.87'
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '1.750'
$ws.Range("E48").Value = '  -1.89%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.235'
$ws.Range("E49").Value = '  -6.08%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05908'
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("B51").Value = 'XinFinNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D51").Value = '0.07240'
$ws.Range("E51").Value = '  +5.71%  '

# Restore default styling on column D now that the text values are in place.
$dRange.Style = "Normal"
